$d = $word.ActiveDocument

function ReplaceInParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $f = $rng.Find
    $f.ClearFormatting()
    $ok = $f.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Host "FAILED to replace in paragraph $paraIndex : $oldText"
    }
    return $ok
}

function ItalicizeFirstInRange($rng, $text) {
    $f = $rng.Find
    $f.ClearFormatting()
    $ok = $f.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $rng.Italic = 1
    } else {
        Write-Host "FAILED to find for italics: $text"
    }
    return $ok
}

# ---------------------------------------------------------------------------
# Entrance Ticket, item 1: "The function of __ ... " ->
#   "Noun: The function of a __ ..." with "function" italicized.
# ---------------------------------------------------------------------------
ReplaceInParagraph 3 "The function of " "Noun: The function of a "
$p3 = $d.Paragraphs(3)
ItalicizeFirstInRange $p3.Range "function"

# ---------------------------------------------------------------------------
# Entrance Ticket, item 2: "If ___ ... doesn't / don't function, ... " ->
#   "Verb: If your ___ ... doesn't function, ... " with "function" italicized
#   and the " / don't" alternative removed.
# ---------------------------------------------------------------------------
ReplaceInParagraph 10 "If ___________________________ " "Verb: If your ___________________________ "
ReplaceInParagraph 10 "/ don’t " ""
$p10 = $d.Paragraphs(10)
ItalicizeFirstInRange $p10.Range "function"

# ---------------------------------------------------------------------------
# Exit Ticket, item 1: italicize "independent variable" and "dependent variable".
# ---------------------------------------------------------------------------
$p18 = $d.Paragraphs(18)
$rng18 = $p18.Range
$ok = ItalicizeFirstInRange $rng18 "independent variable"
$rng18b = $d.Range($rng18.End, $p18.Range.End)
ItalicizeFirstInRange $rng18b "dependent variable"

# ---------------------------------------------------------------------------
# Exit Ticket, item 2: "One way to tell a function from a non-function is ..." ->
#   "Here's my definition of a function: A function is ..." with both
#   "function" instances italicized, and the _GoBack bookmark repositioned
#   around the second "function " (between the colon clause and "is ...").
# ---------------------------------------------------------------------------
ReplaceInParagraph 22 "One way to tell a function from a non-function is …" "Here’s my definition of a function: A function is …"

$p22 = $d.Paragraphs(22)
$rng22 = $p22.Range
ItalicizeFirstInRange $rng22 "function"

$rng22b = $d.Range($rng22.End, $p22.Range.End)
ItalicizeFirstInRange $rng22b "function"

$bmStart = $rng22b.Start
$bmEnd = $rng22b.End + 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmStart, $bmEnd))

# ---------------------------------------------------------------------------
# Exit Ticket, item 3: "surprised me" -> "I wondered about".
# ---------------------------------------------------------------------------
ReplaceInParagraph 26 "surprised me" "I wondered about"
